# Updated cryptos list on Sun Nov 12 15:06:52 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for each coin row, and
# swap the ARBITRUM / InjectiveProtocol rows (43 <-> 44) to their new
# ranking positions with their refreshed figures.
#
# Numeric-looking Price values are entered with a leading apostrophe
# (Excel's normal "store as text" convention) so they stay text cells,
# matching the rest of the column (e.g. "37.208.74") which already can't
# be parsed as plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.208.74'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '2.061.17'
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''249.32'
$ws.Range("E5").Value = '  -1.37%  '
$ws.Range("D6").Value = '''0.667'
$ws.Range("E6").Value = '  -1.38%  '
$ws.Range("D7").Value = '''59.51'
$ws.Range("E7").Value = '  -4.15%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '''0.387'
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").Value = '''0.0789'
$ws.Range("E10").Value = '  -2.35%  '
$ws.Range("D11").Value = '''0.109'
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").Value = '''15.94'
$ws.Range("E12").Value = '  +1.80%  '
$ws.Range("D13").Value = '2.357.97'
$ws.Range("E13").Value = '  -1.21%  '
$ws.Range("D14").Value = '''0.836'
$ws.Range("E14").Value = '  +0.83%  '
$ws.Range("E15").Value = '  +5.97%  '
$ws.Range("D16").Value = '2.053.80'
$ws.Range("E16").Value = '  -1.52%  '
$ws.Range("D17").Value = '''18.15'
$ws.Range("E17").Value = '  +20.72%  '
$ws.Range("D18").Value = '37.216.27'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").Value = '''75.18'
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").Value = '0.0₃0904'
$ws.Range("E20").Value = '  -2.40%  '
$ws.Range("D21").Value = '''5.42'
$ws.Range("E21").Value = '  -0.95%  '
$ws.Range("D22").Value = '''238.25'
$ws.Range("E22").Value = '  -0.93%  '
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").Value = '''2.50'
$ws.Range("E24").Value = '  +3.40%  '
$ws.Range("D25").Value = '''2.20'
$ws.Range("E25").Value = '  +1.97%  '
$ws.Range("D26").Value = '''169.40'
$ws.Range("E26").Value = '  -1.36%  '
$ws.Range("D27").Value = '''9.44'
$ws.Range("E27").Value = '  +1.64%  '
$ws.Range("D28").Value = '''20.17'
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").Value = '''4.86'
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  +2.75%  '
$ws.Range("D32").Value = '''0.0624'
$ws.Range("E32").Value = '  -2.08%  '
$ws.Range("D33").Value = '''4.59'
$ws.Range("E33").Value = '  +2.60%  '
$ws.Range("D34").Value = '''0.0908'
$ws.Range("E34").Value = '  +0.87%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").Value = '''2.31'
$ws.Range("E36").Value = '  -0.76%  '
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("E38").Value = '  -1.15%  '
$ws.Range("E39").Value = '  -3.00%  '
$ws.Range("D40").Value = '''3.12'
$ws.Range("E40").Value = '  +11.06%  '
$ws.Range("D41").Value = '''5.14'
$ws.Range("E41").Value = '  +10.70%  '
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = '''17.57'
$ws.Range("E43").Value = '  -5.49%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '''1.16'
$ws.Range("E44").Value = '  -0.67%  '
$ws.Range("D45").Value = '''96.91'
$ws.Range("E45").Value = '  -2.26%  '
$ws.Range("D46").Value = '''2.49'
$ws.Range("E46").Value = '  -1.96%  '
$ws.Range("D47").Value = '1.290.98'
$ws.Range("E47").Value = '  -1.76%  '
$ws.Range("D48").Value = '''2.92'
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("D49").Value = '''6.92'
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("D50").Value = '2.246.53'
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("D51").Value = '''3.60'
$ws.Range("E51").Value = '  -16.78%  '
